# Auto-generated edit script for violent-crime-ytd workbook update (2023-09-13 data)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 5346
$ws.Range("J3").Value = 5673
$ws.Range("F4").Value = 1362
$ws.Range("J4").Value = 1254
$ws.Range("J5").Value = 443
$ws.Range("J6").Value = 7102
$ws.Range("F7").Value = 17002
$ws.Range("J7").Value = 19818

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J6").Value = 170
$ws.Range("J7").Value = 268

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 379
$ws.Range("J6").Value = 417
$ws.Range("J7").Value = 1250

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("J2").Value = 122
$ws.Range("J3").Value = 155
$ws.Range("J7").Value = 412

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J6").Value = 312
$ws.Range("J7").Value = 913

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J2").Value = 185
$ws.Range("J6").Value = 176
$ws.Range("J7").Value = 617

$ws = $wb.Worksheets.Item("New City")
$ws.Range("J2").Value = 148
$ws.Range("J7").Value = 509

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("J2").Value = 28
$ws.Range("J7").Value = 80

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J2").Value = 159
$ws.Range("J4").Value = 77
$ws.Range("J7").Value = 578
$ws.Range("J8").Value = 1250
$ws.Range("J10").Value = 134
$ws.Range("J11").Value = 308
$ws.Range("J15").Value = 219
$ws.Range("J19").Value = 569
$ws.Range("J20").Value = 414
$ws.Range("J27").Value = 118
$ws.Range("J29").Value = 1107
$ws.Range("J30").Value = 80
$ws.Range("J31").Value = 175
$ws.Range("J32").Value = 31
$ws.Range("J33").Value = 913
$ws.Range("J34").Value = 95
$ws.Range("J37").Value = 617
$ws.Range("J39").Value = 7
$ws.Range("J43").Value = 166
$ws.Range("J46").Value = 67
$ws.Range("J48").Value = 229
$ws.Range("J49").Value = 132
$ws.Range("J53").Value = 268
$ws.Range("J54").Value = 383
$ws.Range("J55").Value = 260
$ws.Range("F63").Value = 130
$ws.Range("J63").Value = 68
$ws.Range("J64").Value = 135
$ws.Range("J65").Value = 509
$ws.Range("J67").Value = 759
$ws.Range("J73").Value = 183
$ws.Range("J77").Value = 157
$ws.Range("J78").Value = 246
$ws.Range("J79").Value = 566
$ws.Range("J83").Value = 412
$ws.Range("J85").Value = 843
$ws.Range("J86").Value = 121
$ws.Range("J88").Value = 218
$ws.Range("J90").Value = 216
$ws.Range("J91").Value = 220
$ws.Range("J92").Value = 60
$ws.Range("J94").Value = 197
$ws.Range("F101").Value = 17002
$ws.Range("J101").Value = 19818

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J6").Value = 48
$ws.Range("J7").Value = 175

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 190
$ws.Range("J6").Value = 199
$ws.Range("J7").Value = 759

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 132

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 73
$ws.Range("J7").Value = 383

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J3").Value = 381
$ws.Range("J7").Value = 1107

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J4").Value = 36
$ws.Range("J7").Value = 229

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J6").Value = 209
$ws.Range("J7").Value = 569

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("J5").Value = 1
$ws.Range("J7").Value = 134

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J6").Value = 68
$ws.Range("J7").Value = 246

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("J6").Value = 127
$ws.Range("J7").Value = 260

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 67

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 220

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 200
$ws.Range("J6").Value = 158
$ws.Range("J7").Value = 566

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 135

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("J2").Value = 114
$ws.Range("J3").Value = 143
$ws.Range("J7").Value = 414

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J2").Value = 180
$ws.Range("J3").Value = 176
$ws.Range("J4").Value = 24
$ws.Range("J7").Value = 578

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 95

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 107
$ws.Range("J7").Value = 197

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 53
$ws.Range("J6").Value = 93
$ws.Range("J7").Value = 219

$ws = $wb.Worksheets.Item("East Village")
$ws.Range("J3").Value = 1

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("J6").Value = 7

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J6").Value = 121
$ws.Range("J7").Value = 308

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("J2").Value = 65
$ws.Range("J7").Value = 183

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("J2").Value = 46
$ws.Range("J3").Value = 39
$ws.Range("J7").Value = 159

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J2").Value = 45
$ws.Range("J7").Value = 218

$ws = $wb.Worksheets.Item("Galewood")
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("J3").Value = 28
$ws.Range("J6").Value = 40
$ws.Range("J7").Value = 118

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("J3").Value = 18
$ws.Range("J7").Value = 121

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 62
$ws.Range("J7").Value = 216

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("J3").Value = 33
$ws.Range("J7").Value = 166

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J2").Value = 219
$ws.Range("J3").Value = 311
$ws.Range("J7").Value = 843

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J2").Value = 57
$ws.Range("J7").Value = 157

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("J5").Value = 2
$ws.Range("J7").Value = 77
